$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update the Date value (Property/Value pair on row 8)
$ws.Range("B8").Value = "2025-05-02T13:48:14+00:00"

# Add two new metadata rows: Source / Target, mirroring the
# ConceptMap source & target StructureDefinition URLs already present
# on the "Mapping Table 0" sheet.
$ws.Range("A15").Value = "Source"
$ws.Range("B15").Value = "https://interop.esante.gouv.fr/ig/fhir/pdsm4dmp/StructureDefinition/DocumentEntry"

$ws.Range("A16").Value = "Target"
$ws.Range("B16").Value = "https://interop.esante.gouv.fr/ig/fhir/pdsm/StructureDefinition/pdsm-comprehensive-document-reference"

# Copy the formatting of the last pre-existing row (14) onto the two
# new rows so they match the rest of the table's style.
$ws.Range("A14:B14").Copy()
$ws.Range("A15:B16").PasteSpecial(-4122)
